# Correzione errori di battitura (typo fixes) in "data dictionary comune.xlsx"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) E4: "Stazione stadale, cellula stradale" -> "Stazione stradale, cellula stradale"
$ws.Range("E4").Value = "Stazione stradale, cellula stradale"

# 2) B11: "...oppure anche da quest’ultime verso l’utente." -> "...oppure anche da queste ultime verso l’utente."
$ws.Range("B11").Value = "Segnalazione dell’avvenimento di un dato evento stradale riguardante il traffico. Può essere indirizzata da un sottosistema verso il sistema centrale o dal sistema centrale verso le applicazioni mobili, oppure anche da queste ultime verso l’utente."

# 3) K14: "Notifica, dati di traffco, mappa, conteggio, posizione," -> "Notifica, dati di traffico, mappa, conteggio, posizione,"
$ws.Range("K14").Value = "Notifica, dati di traffico, mappa, conteggio, posizione,"

# 4) B15: "...l’app moble)" -> "...l’app mobile)"
$ws.Range("B15").Value = "Automobile su cui è installata una centralina auto (oppure unità contata dalla centralina stradale oppure mezzo di trasporto dell’utente che comunica con l’app mobile)"
